$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.138.26"
$ws.Range("E2").Value = "  +0.95%  "

$ws.Range("D3").Value = "2.049.96"
$ws.Range("E3").Value = "  -3.08%  "

$ws.Range("E4").Value = "  +0.06%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "249.28"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.85%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.655"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.34%  "

$ws.Range("E7").Value = "  -0.10%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "55.56"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +16.13%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "61.88"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.42%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.378"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.01%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0756"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.16%  "

$ws.Range("E12").Value = "  +5.70%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "15.12"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +4.79%  "

$ws.Range("D14").Value = "2.347.43"
$ws.Range("E14").Value = "  -3.21%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.825"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.95%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.25"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.37%  "

$ws.Range("D17").Value = "2.049.49"
$ws.Range("E17").Value = "  -3.18%  "

$ws.Range("D18").Value = "37.016.64"
$ws.Range("E18").Value = "  +0.44%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "72.45"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.01%  "

$ws.Range("D20").Value = "0.0₃0884"
$ws.Range("E20").Value = "  +4.85%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.33"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +6.48%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.28"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.59%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "237.73"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.79%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.42"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.09%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "170.46"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.93%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.15"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.10%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "20.29"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -5.54%  "

$ws.Range("E29").Value = "  -1.07%  "

$ws.Range("E30").Value = "  -0.58%  "

$ws.Range("E31").Value = "  +1.99%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0626"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +4.11%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.05"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +13.38%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.36"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.84%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.28"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.97%  "

$ws.Range("B37").Value = "Gas"
$ws.Range("C37").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "19.15"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -23.59%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0849"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -10.45%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.79"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.96%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.110"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +31.45%  "

$ws.Range("E41").Value = "  -1.20%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "18.27"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +12.50%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0225"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.27%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.14"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.92%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "97.23"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.06%  "

$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.38"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +62.94%  "

$ws.Range("E47").Value = "  -1.35%  "

$ws.Range("D48").Value = "1.304.09"
$ws.Range("E48").Value = "  -4.29%  "

$ws.Range("E49").Value = "  +2.93%  "

$ws.Range("E50").Value = "  +2.87%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "6.87"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -3.55%  "
